$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "HFT setz." rows (bottom-up so row numbers stay valid)
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(7).Delete()

# Append the six new "prefabricated parts" rows (rows 12-17) for the
# "3_5_T1_SP_GRU_OG1_4101_03_F-P-001 - Decke_Kein BA" task group
$newRows = @(
  @("3_5_T1_SP_GRU_OG1_4101_03_F-P-001 - Decke_Kein BA_Beton_BET", "BET", 44950, 44951),
  @("3_5_T1_SP_GRU_OG1_4101_03_F-P-001 - Decke_Kein BA_Beton_BEW", "BEW", 44949, 44950),
  @("3_5_T1_SP_GRU_OG1_4101_03_F-P-001 - Decke_Kein BA_Beton_SCH", "SCH", 44946, 44947),
  @("3_5_T1_SP_GRU_OG1_4101_03_F-P-001 - Decke_Kein BA_HFT_BET", "BET", 44950, 44951),
  @("3_5_T1_SP_GRU_OG1_4101_03_F-P-001 - Decke_Kein BA_HFT_BEW", "BEW", 44949, 44950),
  @("3_5_T1_SP_GRU_OG1_4101_03_F-P-001 - Decke_Kein BA_HFT_HFT", "HFT", 44946, 44947)
)

$startRow = 12
for ($i = 0; $i -lt $newRows.Count; $i++) {
  $row = $startRow + $i
  $values = $newRows[$i]

  $ws.Cells.Item($row, 1).Value = $values[0]
  $ws.Cells.Item($row, 2).Value = $values[1]

  $ws.Cells.Item($row, 3).Value = $values[2]
  $ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD HH:MM:SS"

  $ws.Cells.Item($row, 4).Value = $values[3]
  $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
